$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(39400, 2007, 3.145939949069287, 2008, 2.455972543253826),
    @(39765, 2008, 1.769627576887389, 2009, 1.087227286828241),
    @(40130, 2009, -4.774178217057779, 2010, -0.8523446516643385),
    @(40494, 2010, 1.97975191822708, 2011, 1.39052144387346),
    @(40862, 2011, 3.452886745653183, 2012, 2.798447799311043),
    @(41228, 2012, 1.239479831392853, 2013, 2.047133666472267),
    @(41592, 2013, 0.2379616621361214, 2014, 1.062273708599726),
    @(41957, 2014, 1.51977456621637, 2015, 0.6176326357196116),
    @(42321, 2015, 1.470039379455756, 2016, 1.577608035818301),
    @(42689, 2016, 1.638797242243251, 2017, 1.369334405341616),
    @(43053, 2017, 2.161565493242668, 2018, 2.486299099038347),
    @(43418, 2018, 2.214251681313772, 2019, 1.362030665126834),
    @(43783, 2019, 0.6066442151010376, 2020, 0.7878236429522678),
    @(44159, 2020, -4.207901339433196, 2021, -0.4850133725290084),
    @(44525, 2021, 1.099928004397532, 2022, 1.102608990163567),
    @(44890, 2022, 2.310042359896225, 2023, 2.086256540666986),
    @(45254, 2023, 0.0464415346324687, 2024, 0.3491198177708599),
    @(45618, 2024, -0.3101476031197148, 2025, 0.2126457877301924)
)

# Extend the A-column date style down to the new row 19 before writing values,
# mirroring the style used by the existing date cells (A2:A18).
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
